$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''70.353.81'
$ws.Range("E2").Value = '  +4.70%  '
$ws.Range("D3").Value = '''3.609.98'
$ws.Range("E3").Value = '  +4.63%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '''585.92'
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("D6").Value = '''192.46'
$ws.Range("E6").Value = '  +2.48%  '
$ws.Range("D7").Value = '''0.637'
$ws.Range("E7").Value = '  +1.10%  '
$ws.Range("D8").Value = '''3.604.20'
$ws.Range("E8").Value = '  +4.67%  '
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("D10").Value = '''0.181'
$ws.Range("E10").Value = '  +5.26%  '
$ws.Range("D11").Value = '''0.667'
$ws.Range("E11").Value = '  +3.89%  '
$ws.Range("D12").Value = '''57.56'
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D13").Value = '''0.0000304'
$ws.Range("E13").Value = '  +10.19%  '
$ws.Range("D14").Value = '''9.81'
$ws.Range("E14").Value = '  +3.97%  '
$ws.Range("D15").Value = '''4.213.70'
$ws.Range("E15").Value = '  +5.50%  '
$ws.Range("D16").Value = '''20.15'
$ws.Range("E16").Value = '  +5.99%  '
$ws.Range("D17").Value = '''3.631.24'
$ws.Range("E17").Value = '  +5.64%  '
$ws.Range("D18").Value = '''70.576.88'
$ws.Range("E18").Value = '  +5.19%  '
$ws.Range("D19").Value = '''12.59'
$ws.Range("E19").Value = '  +4.43%  '
$ws.Range("E20").Value = '  +2.67%  '
$ws.Range("D21").Value = '''1.05'
$ws.Range("E21").Value = '  +3.31%  '
$ws.Range("D22").Value = '''483.15'
$ws.Range("E22").Value = '  -1.34%  '
$ws.Range("D23").Value = '''19.32'
$ws.Range("E23").Value = '  +13.65%  '
$ws.Range("D24").Value = '''5.08'
$ws.Range("E24").Value = '  -10.02%  '
$ws.Range("D25").Value = '''4.44'
$ws.Range("E25").Value = '  +2.65%  '
$ws.Range("D26").Value = '''90.12'
$ws.Range("E26").Value = '  +0.60%  '
$ws.Range("D27").Value = '''3.12'
$ws.Range("E27").Value = '  +5.01%  '
$ws.Range("D28").Value = '''11.31'
$ws.Range("E28").Value = '  +3.32%  '
$ws.Range("D29").Value = '''9.42'
$ws.Range("E29").Value = '  +4.49%  '
$ws.Range("D30").Value = '''8.00'
$ws.Range("E30").Value = '  +9.00%  '
$ws.Range("D31").Value = '''32.43'
$ws.Range("E31").Value = '  +3.80%  '
$ws.Range("D32").Value = '''0.120'
$ws.Range("E32").Value = '  +7.08%  '
$ws.Range("D33").Value = '''12.20'
$ws.Range("E33").Value = '  +3.41%  '
$ws.Range("D34").Value = '''66.55'
$ws.Range("E34").Value = '  +2.61%  '
$ws.Range("D35").Value = '''610.69'
$ws.Range("E35").Value = '  +0.94%  '
$ws.Range("D36").Value = '''40.08'
$ws.Range("E36").Value = '  +8.38%  '
$ws.Range("D37").Value = '''0.0₃0830'
$ws.Range("E37").Value = '  +6.07%  '
$ws.Range("D38").Value = '''0.406'
$ws.Range("E38").Value = '  +5.16%  '
$ws.Range("D39").Value = '''0.147'
$ws.Range("E39").Value = '  +1.21%  '
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").Value = '''2.96'
$ws.Range("E41").Value = '  +15.01%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''3.56'
$ws.Range("E42").Value = '  +3.03%  '
$ws.Range("D43").Value = '''3.311.11'
$ws.Range("E43").Value = '  +3.77%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = '''3.16'
$ws.Range("E44").Value = '  +18.75%  '
$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D45").Value = '''3.12'
$ws.Range("E45").Value = '  +8.04%  '
$ws.Range("D46").Value = '''0.0453'
$ws.Range("E46").Value = '  +5.68%  '
$ws.Range("D47").Value = '''9.66'
$ws.Range("E47").Value = '  +11.88%  '
$ws.Range("D48").Value = '''3.40'
$ws.Range("E48").Value = '  +5.39%  '
$ws.Range("D49").Value = '''0.138'
$ws.Range("E49").Value = '  +2.08%  '
$ws.Range("D50").Value = '''0.999'
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("D51").Value = '''3.19'
$ws.Range("E51").Value = '  +0.65%  '
